$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update classification report values (new t-test based metrics)
$ws.Range("B2").Value = 0.7828668363019508
$ws.Range("C2").Value = 0.7984429065743944
$ws.Range("D2").Value = 0.7905781584582441
$ws.Range("E2").Value = 1156

$ws.Range("B3").Value = 0.8717948717948718
$ws.Range("C3").Value = 0.785824345146379
$ws.Range("D3").Value = 0.826580226904376
$ws.Range("E3").Value = 649

$ws.Range("B4").Value = 0.7551282051282051
$ws.Range("C4").Value = 0.7474619289340102
$ws.Range("D4").Value = 0.7512755102040817
$ws.Range("E4").Value = 788

$ws.Range("B5").Value = 0.4974747474747475
$ws.Range("C5").Value = 0.5677233429394812
$ws.Range("D5").Value = 0.5302826379542396
$ws.Range("E5").Value = 347

$ws.Range("B6").Value = 0.7547619047619047
$ws.Range("C6").Value = 0.7547619047619047
$ws.Range("D6").Value = 0.7547619047619047
$ws.Range("E6").Value = 0.7547619047619047

$ws.Range("B7").Value = 0.7268161651749437
$ws.Range("C7").Value = 0.7248631308985662
$ws.Range("D7").Value = 0.7246791333802354

$ws.Range("B8").Value = 0.761378808698874
$ws.Range("C8").Value = 0.7547619047619047
$ws.Range("D8").Value = 0.7572694203570095
